$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Status(Summary)")

# --- Fix a typo in the existing last comment (row 20 / column D) ---
$d20 = $ws.Range("D20").Value2
$d20New = $d20 -replace '7\) Present the findings and screenshots during syncup', '7) Presented the findings and screenshots during syncup'
$ws.Range("D20").Value = $d20New

# --- Add a new status row to the table (Table6) ---
$lo = $ws.ListObjects.Item(1)
$newListRow = $lo.ListRows.Add()
$newRowIndex = $newListRow.Range.Row

# Copy formatting from the previous last data row so the new row matches exactly
# (reuses the existing cell styles instead of synthesizing new ones)
$prevRow = $newRowIndex - 1
$ws.Range("A$prevRow`:D$prevRow").Copy()
$ws.Range("A$newRowIndex`:D$newRowIndex").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

$ws.Range("A$newRowIndex").Value = 44007
$ws.Range("B$newRowIndex").Value = "Implementation"
$ws.Range("C$newRowIndex").Value = "App state management "

$comment = @'
1) Input data captured from textfields, formfields, radio buttons and checkboxes moved into app state
2) Bug fixes during migration of elements from data stored in internal state to app state `Page Provider` 
3) Started a simulation of rerendering and displaying survey items with preset (recorded) values on response change of each survey item 
'@
$ws.Range("D$newRowIndex").Value = $comment

$ws.Rows.Item($newRowIndex).RowHeight = 86

# --- Update selection to follow the newly added last cell ---
$ws.Range("D$newRowIndex").Select() | Out-Null
